# Update the cryptocurrency price table (rows 2-51) on the active sheet
# to reflect the latest scrape: refreshed Price (D) / Volume(1h) (E) values,
# and the FraxShare / TheSandbox rows (39-40) swapping places.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '29.443.92'
$ws.Range("E2").Value = '  -0.20%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '1.900.54'
$ws.Range("E3").Value = '  -0.60%  '

# Row 4: TetherUSD
$ws.Range("D4").Value = '''1.005'
$ws.Range("E4").Value = '  +0.38%  '

# Row 5: BNB
$ws.Range("D5").Value = '''326.26'
$ws.Range("E5").Value = '  -2.01%  '

# Row 6: USDC
$ws.Range("D6").Value = '''1.005'
$ws.Range("E6").Value = '  +0.38%  '

# Row 7: XRP
$ws.Range("D7").Value = '''0.4793'
$ws.Range("E7").Value = '  +2.58%  '

# Row 8: Cardano
$ws.Range("D8").Value = '''0.4059'
$ws.Range("E8").Value = '  -1.06%  '

# Row 9: Dogecoin
$ws.Range("D9").Value = '''0.08073'
$ws.Range("E9").Value = '  +0.50%  '

# Row 10: Polygon
$ws.Range("D10").Value = '''1.002'
$ws.Range("E10").Value = '  -1.07%  '

# Row 11: Solana
$ws.Range("D11").Value = '''23.40'
$ws.Range("E11").Value = '  +4.70%  '

# Row 12: WrappedEther
$ws.Range("D12").Value = '1.909.05'
$ws.Range("E12").Value = '  -0.58%  '

# Row 13: Polkadot
$ws.Range("D13").Value = '''5.960'
$ws.Range("E13").Value = '  -0.24%  '

# Row 14: Chainlink
$ws.Range("D14").Value = '''7.083'
$ws.Range("E14").Value = '  -1.32%  '

# Row 15: Litecoin
$ws.Range("D15").Value = '''90.18'
$ws.Range("E15").Value = '  +0.35%  '

# Row 16: BinanceUSD
$ws.Range("E16").Value = '  +0.46%  '

# Row 17: TRON
$ws.Range("D17").Value = '''0.06723'
$ws.Range("E17").Value = '  +2.11%  '

# Row 18: ShibaInu
$ws.Range("D18").Value = '''0.00001032'
$ws.Range("E18").Value = '  -0.10%  '

# Row 19: Avalanche
$ws.Range("D19").Value = '''17.62'
$ws.Range("E19").Value = '  -0.84%  '

# Row 20: Dai
$ws.Range("D20").Value = '''1.005'
$ws.Range("E20").Value = '  +0.34%  '

# Row 21: WrappedBTC
$ws.Range("D21").Value = '29.465.12'
$ws.Range("E21").Value = '  +0.00%  '

# Row 22: Uniswap
$ws.Range("D22").Value = '''5.546'
$ws.Range("E22").Value = '  -0.42%  '

# Row 23: Cosmos
$ws.Range("E23").Value = '  +2.39%  '

# Row 24: Toncoin
$ws.Range("D24").Value = '''2.158'
$ws.Range("E24").Value = '  -2.52%  '

# Row 25: WrappedliquidstakedEther2.0
$ws.Range("D25").Value = '2.188.68'
$ws.Range("E25").Value = '  +2.01%  '

# Row 26: Monero
$ws.Range("D26").Value = '''154.14'
$ws.Range("E26").Value = '  -0.40%  '

# Row 27: EthereumClassic
$ws.Range("D27").Value = '''19.88'
$ws.Range("E27").Value = '  -0.06%  '

# Row 28: InternetComputer(DFINITY)
$ws.Range("D28").Value = '''6.095'
$ws.Range("E28").Value = '  +5.62%  '

# Row 29: LidoDAOToken
$ws.Range("D29").Value = '''2.091'
$ws.Range("E29").Value = '  -2.39%  '

# Row 30: BitcoinCash
$ws.Range("D30").Value = '''118.50'
$ws.Range("E30").Value = '  +0.95%  '

# Row 31: ImmutableX
$ws.Range("D31").Value = '''1.032'
$ws.Range("E31").Value = '  -3.05%  '

# Row 32: Stellar
$ws.Range("D32").Value = '''0.09487'
$ws.Range("E32").Value = '  +0.24%  '

# Row 33: Filecoin
$ws.Range("D33").Value = '''5.510'
$ws.Range("E33").Value = '  +1.97%  '

# Row 34: HuobiToken
$ws.Range("E34").Value = '  -0.70%  '

# Row 35: ARBITRUM
$ws.Range("D35").Value = '''1.389'
$ws.Range("E35").Value = '  -2.62%  '

# Row 36: Hedera
$ws.Range("D36").Value = '''0.06080'
$ws.Range("E36").Value = '  -0.59%  '

# Row 37: VeChain
$ws.Range("D37").Value = '''0.02252'
$ws.Range("E37").Value = '  -0.41%  '

# Row 38: TrustWalletToken
$ws.Range("D38").Value = '''1.171'
$ws.Range("E38").Value = '  -0.76%  '

# Row 39: FraxShare
$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").Value = '''7.945'
$ws.Range("E39").Value = '  -5.60%  '

# Row 40: TheSandbox
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = '''0.5883'
$ws.Range("E40").Value = '  -0.08%  '

# Row 41: Algorand
$ws.Range("D41").Value = '''0.1844'
$ws.Range("E41").Value = '  +0.11%  '

# Row 42: Aptos
$ws.Range("D42").Value = '''10.24'
$ws.Range("E42").Value = '  +0.29%  '

# Row 43: WEMIXToken
$ws.Range("D43").Value = '''1.292'

# Row 44: Cronos
$ws.Range("D44").Value = '''0.07817'
$ws.Range("E44").Value = '  +4.12%  '

# Row 45: RenderToken
$ws.Range("D45").Value = '''2.390'
$ws.Range("E45").Value = '  +1.48%  '

# Row 46: EnergySwap
$ws.Range("D46").Value = '''12.22'
$ws.Range("E46").Value = '  +0.29%  '

# Row 47: Decentraland
$ws.Range("D47").Value = '''0.5533'
$ws.Range("E47").Value = '  -0.59%  '

# Row 48: NEARProtocol
$ws.Range("D48").Value = '''1.922'
$ws.Range("E48").Value = '  -0.25%  '

# Row 49: Quant
$ws.Range("D49").Value = '''114.27'
$ws.Range("E49").Value = '  +0.90%  '

# Row 50: Aave
$ws.Range("D50").Value = '''72.36'
$ws.Range("E50").Value = '  +1.07%  '

# Row 51: WOONetwork
$ws.Range("D51").Value = '''0.2935'
$ws.Range("E51").Value = '  -1.02%  '
